$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (e.g. "60.521.97"
# using dots as separators, or values with significant trailing zeros like
# "0.120"). Force text number format on each target cell before assigning so
# Excel does not reinterpret the literal as a number and mangle it.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.521.97'
$ws.Range("E2").Value = '  -3.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.313.78'
$ws.Range("E3").Value = '  -3.63%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.84'
$ws.Range("E5").Value = '  -3.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.69'
$ws.Range("E6").Value = '  -3.59%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.313.62'
$ws.Range("E8").Value = '  -3.62%  '

$ws.Range("E9").Value = '  -1.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.85'
$ws.Range("E10").Value = '  -2.15%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.120'
$ws.Range("E11").Value = '  -3.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.411'
$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.869.17'

$ws.Range("E14").Value = '  +0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.25'
$ws.Range("E15").Value = '  -3.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.293.62'
$ws.Range("E16").Value = '  -4.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000167'
$ws.Range("E17").Value = '  -3.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.410.72'
$ws.Range("E18").Value = '  -3.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("E19").Value = '  -3.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.45'
$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.65'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '376.32'
$ws.Range("E22").Value = '  -2.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.98'
$ws.Range("E23").Value = '  -1.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.548'
$ws.Range("E24").Value = '  -3.45%  '

$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.436.65'

$ws.Range("E27").Value = '  -7.82%  '

$ws.Range("E28").Value = '  -6.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.30'
$ws.Range("E30").Value = '  -4.76%  '

$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.71'
$ws.Range("E32").Value = '  -3.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.05'
$ws.Range("E33").Value = '  -3.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.63'
$ws.Range("E34").Value = '  -2.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("E35").Value = '  -5.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.25'
$ws.Range("E36").Value = '  -3.70%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.56'
$ws.Range("E37").Value = '  -4.95%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '166.32'
$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.77'
$ws.Range("E39").Value = '  -2.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.74'
$ws.Range("E40").Value = '  -15.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0746'
$ws.Range("E42").Value = '  -5.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.06'
$ws.Range("E43").Value = '  -0.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.753'
$ws.Range("E44").Value = '  -3.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.22'
$ws.Range("E45").Value = '  -3.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.61'
$ws.Range("E46").Value = '  -5.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.13'
$ws.Range("E47").Value = '  -3.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.357.60'
$ws.Range("E48").Value = '  -7.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.57'
$ws.Range("E50").Value = '  -5.24%  '

$ws.Range("E51").Value = '  -3.77%  '
